$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price / volume(1h) table cells to match the latest scraped data.
# Columns D (Price) and E (Volume(1h)) store numeric-looking text (e.g. "311.03", "1.00%")
# as literal strings, so each such cell is forced to Text format ("@") before the
# assignment to prevent Excel from re-interpreting it as a number/percentage and
# silently dropping significant trailing zeros or switching to scientific notation.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '311.03'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.00%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '38.13'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.13%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.116'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.36%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07917'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.49%'
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.899'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-3.53%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.238'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.13%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.842'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-8.20%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9272'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.44%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1200'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-7.76%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1899'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.67%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09270'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '5.23%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03378'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-1.55%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09594'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.41%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001363'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.99%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005820'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.77%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.526'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.05%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.398'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.20%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3447'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.39%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.281'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '5.71%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1281'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.60%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '4.02%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '179.84%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.04372'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.32%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001249'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '2.55%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004274'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-7.39%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001300'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-63.80%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02102'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-9.44%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05075'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.83%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007661'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.93%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009103'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-7.91%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1351'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-0.46%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002021'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-3.19%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.008635'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '7.89%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006672'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '1.79%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.10%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002901'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-3.16%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001196'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.34%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.10%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.10%'
